# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Row 26/27 also swap: PEPE moves above Binance-PegBSC-USD in the ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.126.75'
$ws.Range("E2").Value = '  +3.59%  '

$ws.Range("D3").Value = '3.066.96'
$ws.Range("E3").Value = '  +6.52%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''513.89'
$ws.Range("E5").Value = '  +5.42%  '

$ws.Range("D6").Value = '''140.35'
$ws.Range("E6").Value = '  +7.20%  '

$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +4.45%  '

$ws.Range("E9").Value = '  +1.29%  '

$ws.Range("E10").Value = '  +5.51%  '

$ws.Range("D11").Value = '''0.370'
$ws.Range("E11").Value = '  +7.61%  '

$ws.Range("D12").Value = '3.583.95'
$ws.Range("E12").Value = '  +6.26%  '

$ws.Range("E13").Value = '  +3.24%  '

$ws.Range("D14").Value = '''25.23'
$ws.Range("E14").Value = '  +0.15%  '

$ws.Range("E15").Value = '  +5.73%  '

$ws.Range("D16").Value = '57.213.44'
$ws.Range("E16").Value = '  +3.79%  '

$ws.Range("D17").Value = '3.066.84'
$ws.Range("E17").Value = '  +6.49%  '

$ws.Range("D18").Value = '''5.94'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").Value = '''13.04'
$ws.Range("E19").Value = '  +6.22%  '

$ws.Range("D20").Value = '''8.14'
$ws.Range("E20").Value = '  +8.16%  '

$ws.Range("D21").Value = '''335.44'
$ws.Range("E21").Value = '  +8.08%  '

$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("E23").Value = '  +6.07%  '

$ws.Range("D24").Value = '''65.23'
$ws.Range("E24").Value = '  +5.79%  '

$ws.Range("D25").Value = '''0.168'
$ws.Range("E25").Value = '  +6.48%  '

$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = '0.0₃0953'
$ws.Range("E26").Value = '  +15.47%  '

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.31%  '

$ws.Range("D28").Value = '''6.43'
$ws.Range("E28").Value = '  +2.88%  '

$ws.Range("D29").Value = '''6.97'
$ws.Range("E29").Value = '  +1.46%  '

$ws.Range("E30").Value = '  +5.07%  '

$ws.Range("D31").Value = '''20.75'
$ws.Range("E31").Value = '  +6.41%  '

$ws.Range("E32").Value = '  +7.22%  '

$ws.Range("D33").Value = '''154.62'
$ws.Range("E33").Value = '  +4.27%  '

$ws.Range("E34").Value = '  +4.95%  '

$ws.Range("D35").Value = '''5.85'
$ws.Range("E35").Value = '  +6.50%  '

$ws.Range("D36").Value = '''26.47'
$ws.Range("E36").Value = '  +9.22%  '

$ws.Range("D37").Value = '''1.23'
$ws.Range("E37").Value = '  +5.72%  '

$ws.Range("D38").Value = '''0.0669'
$ws.Range("E38").Value = '  +4.12%  '

$ws.Range("D39").Value = '3.103.24'
$ws.Range("E39").Value = '  +6.64%  '

$ws.Range("D40").Value = '''36.93'
$ws.Range("E40").Value = '  +2.39%  '

$ws.Range("D41").Value = '''0.669'
$ws.Range("E41").Value = '  +6.98%  '

$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("E43").Value = '  +6.07%  '

$ws.Range("D44").Value = '2.232.90'

$ws.Range("E45").Value = '  +10.76%  '

$ws.Range("D47").Value = '''0.938'
$ws.Range("E47").Value = '  +4.55%  '

$ws.Range("D48").Value = '''19.88'
$ws.Range("E48").Value = '  +8.39%  '

$ws.Range("D49").Value = '''5.84'
$ws.Range("E49").Value = '  +0.95%  '

$ws.Range("D50").Value = '''0.0864'
$ws.Range("E50").Value = '  +3.41%  '

$ws.Range("D51").Value = '''0.181'
$ws.Range("E51").Value = '  +6.22%  '
